# RawData_long.xlsx edit
# - Column A ("# Spectrum#"): each data row's spectrum index is decremented by one
#   (row 2 -> 1, row 3 -> 2, ... row 31 -> 30)
# - Column Z ("TSP"): every data value becomes 9 (was 1)
# - Selection moves from V1 (whole column) to Z2:Z31, active cell Z2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1      # column A
    $ws.Cells.Item($r, 26).Value = 9          # column Z
}

[void]$ws.Range("Z2:Z31").Select()
